# Append the newest wallet-info rows (new JSON format -> more snapshot rows)
# to the token price-history sheet: 4 new (Date, USDValue) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 13; Date = "2024-11-15"; Value = 0.02081 },
    @{ Row = 14; Date = "2024-03-16"; Value = 0.01633 },
    @{ Row = 15; Date = "2024-03-17"; Value = 0.01565 },
    @{ Row = 16; Date = "2024-08-20"; Value = 0.01106 }
)

$lastRow = $newRows[$newRows.Count - 1].Row

# Format column A as Text for the new rows first so the yyyy-mm-dd strings are
# stored verbatim (as shared-string text, like the existing Date column)
# instead of being auto-converted into date serial numbers.
$dateCol = $ws.Range("A13:A$lastRow")
$dateCol.NumberFormat = "@"

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Value
}

# Drop back to the default (unstyled) cell format so these new cells match
# the plain look of the pre-existing data rows.
$dateCol.Style = "Normal"
